$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'76.125.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'2.919.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.57%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'203.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.77%  "
$ws.Range("D6").Value = "'597.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("D10").Value = "'2.918.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.59%  "
$ws.Range("D11").Value = "'0.431"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +16.27%  "
$ws.Range("D12").Value = "'0.161"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "'3.458.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.57%  "
$ws.Range("D15").Value = "'76.050.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "'28.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.66%  "
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").Value = "'2.917.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.47%  "
$ws.Range("D19").Value = "'12.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.19%  "
$ws.Range("D20").Value = "'8.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.25%  "
$ws.Range("D21").Value = "'372.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("E22").Value = "  +1.81%  "
$ws.Range("E23").Value = "  +5.45%  "
$ws.Range("D24").Value = "'71.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Value = "'3.075.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.72%  "
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("D28").Value = "'9.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("E29").Value = "  +3.69%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("D32").Value = "'500.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.28%  "
$ws.Range("D33").Value = "'7.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  +2.63%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").Value = "'165.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").Value = "'20.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.66%  "
$ws.Range("D38").Value = "'0.110"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +28.74%  "
$ws.Range("D39").Value = "'19.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("E40").Value = "  -4.86%  "
$ws.Range("D41").Value = "'0.365"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.35%  "
$ws.Range("D42").Value = "'182.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.17%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "'5.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "'39.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("D48").Value = "'2.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("D49").Value = "'0.572"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.09%  "
$ws.Range("D50").Value = "'3.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").Value = "'22.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.33%  "
